# Commit: "changed conclusion to results, put the numbers there"
#
# What actually happened in this edit:
#   1) The "Related Work" slide (originally the LAST slide, position 8) was
#      moved up to right after "Our Approaches", i.e. to position 4 - so it
#      now comes right before the "SRPT" slide.
#   2) The "SRPT" slide lost its last two bullets ("Under 0.8 load..." and
#      "Helps small frame sizes...").
#   3) The "Frame Counter" slide gained a new second bullet about
#      theoretical improvement.
#   4) The "Dynamic Backoff" slide is untouched content-wise (only moved).
#   5) The old "Conclusion" slide (now the last slide again, position 8)
#      was renamed to "Results" and its body was rewritten to include the
#      two bullets removed from the "SRPT" slide (the "numbers").

$p = $ppt.ActivePresentation

# --- 1) Move "Related Work" (slide 8) to slide position 4 -----------------
$relatedWork = $p.Slides.Item(8)
$relatedWork.MoveTo(4)

# After the move the slide order is:
#   1 WiBi (title)           - unchanged
#   2 Intro                  - unchanged
#   3 Our Approaches         - unchanged
#   4 Related Work           - moved here, content unchanged
#   5 SRPT                   - bullets trimmed
#   6 Frame Counter          - new bullet added
#   7 Dynamic Backoff        - unchanged, just shifted down
#   8 Conclusion -> Results  - retitled + rewritten body

# --- 2) SRPT slide (now #5): drop the last two bullets --------------------
$srpt = $p.Slides.Item(5)
$srptBody = $srpt.Shapes.Item(2).TextFrame.TextRange
$srptBody.Text = "Based on Shortest Remaining Processing Time for Processor Scheduling.`r" + `
    "The maximum for the Random Back off is based on the number of Remaining Frames"

# --- 3) Frame Counter slide (now #6): add the theoretical-improvement note -
$frameCounter = $p.Slides.Item(6)
$frameCounterBody = $frameCounter.Shapes.Item(2).TextFrame.TextRange
$frameCounterBody.Text = "The Maximum random back off is based on the current Frame number being sent`r" + `
    "Theoretically, improvement will be similar to that of SRPT, but can handle unknown stream sizes"

# --- 4) Dynamic Backoff slide (now #7): no content change -----------------
# (left as-is; only its slide position shifted because of the move above)

# --- 5) Conclusion -> Results slide (now #8) -------------------------------
$results = $p.Slides.Item(8)
$results.Shapes.Item(1).TextFrame.TextRange.Text = "Results"

$resultsBody = $results.Shapes.Item(2).TextFrame.TextRange
$resultsBody.Text = "It is indeed possible to increase the throughput for 802.11 using a biased MAC protocol using a model based off of SRPT. `r" + `
    "Under 0.8 load, 80% of node traffic improves by factor of 10 `r" + `
    "Helps small frame sizes substantially, only slightly penalizes large frame sizes. `r`r"
